# "added heat recycling for lime and o2"
#
# Renames the two auxiliary-process labels used by the lime and oxygen
# (O2) heat-recovery connections so they read as regular words instead
# of underscore_joined identifiers:
#   aux_lime_kiln      -> aux_lime kiln
#   aux_air_seperation -> aux_air separation
#
# These labels are shared strings referenced from the "ref connections"
# sheet (rows 3/4 - the normal process rows - and rows 17/18 - the new
# heat-recovery rows). Resetting the cell style to Normal before writing
# the value + Text number format matches the target formatting (plain
# text format, default font) instead of the previous custom-font style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ref connections")

$ws.Range("H3").Style = "Normal"
$ws.Range("H3").Value = "aux_lime kiln"
$ws.Range("H3").NumberFormat = "@"

$ws.Range("H4").Style = "Normal"
$ws.Range("H4").Value = "aux_air separation"
$ws.Range("H4").NumberFormat = "@"

$ws.Range("H17").Style = "Normal"
$ws.Range("H17").Value = "aux_air separation"
$ws.Range("H17").NumberFormat = "@"

$ws.Range("H18").Style = "Normal"
$ws.Range("H18").Value = "aux_lime kiln"
$ws.Range("H18").NumberFormat = "@"

# Cursor/selection positions moved in the authored edit.
$wsChains = $wb.Worksheets.Item("ref chains")
$wsChains.Activate()
$wsChains.Range("E5").Select()

$ws.Activate()
$ws.Range("G25").Select()
